$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ranking 1): score, document name, and content all change
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '0.101'
$ws.Range("C2").Value = 'Noticia N° 09.txt'
$ws.Range("D2").Value = 'título: aporte aporte dólar dólar soja, soja, banco banco central central compró compró solo solo usd usd millones millones mercado mercado resumen: resumen: entidad entidad monetaria monetaria encadenó encadenó ruedas ruedas operativas operativas compras compras netas, netas, margen margen achicó achicó primer primer día día agregado agregado tipo tipo cambio cambio especial especial exportaciones exportaciones sojeras sojeras contenido: contenido: jun, jun, sesión sesión usd usd millones millones segmento segmento contado contado spot, spot, liquidaciones liquidaciones programa programa incremento incremento exportador exportador economías economías regionales, regionales, dólar, dólar, aportaron aportaron usd usd millones, millones, volumen volumen aportado aportado ventas ventas soja soja derivados, derivados, gozaron gozaron estándar estándar cambiario cambiario semana semana pasada. pasada. bcra bcra concluyó concluyó intervención intervención cambiaria cambiaria día día saldo saldo neto neto comprador comprador solo solo millones millones dólares. dólares. operaciones operaciones dólar dólar agro agro registran registran tercera tercera etapa etapa ingresos ingresos usd usd millones millones abril abril parte. parte. asimismo, asimismo, empezó empezó balance balance bcra bcra intervención intervención cambiaria cambiaria negativo negativo usd usd millones, millones, mientras mientras junio junio sostiene sostiene im im saldo saldo comprador comprador millones millones dólares. dólares. bcra bcra aceleró aceleró resguardo resguardo escasas escasas reservas, reservas, medida medida afecta afecta provincias provincias municipios, municipios, momentos momentos ministerio ministerio economía economía planea planea nuevo nuevo canje canje voluntario voluntario bonos bonos pesos pesos objetivo objetivo despejar despejar vencimientos. vencimientos. ministro ministro economía, economía, sergio sergio massa, massa, viajará viajará washington washington cerca cerca junio junio procura procura cerrar cerrar readecuación readecuación acuerdo acuerdo fondo fondo monetario monetario internacional internacional fmi, fmi, permitiría permitiría argentina argentina obtener obtener adelanto adelanto desembolsos desembolsos nuevas nuevas metas metas cumplir, cumplir, luego luego impacto impacto sequía sequía exportaciones exportaciones sector sector agro, agro, “las “las negociaciones negociaciones argentina argentina funcionarios funcionarios técnicos técnicos fondo fondo avanzan avanzan hace hace casi casi dos dos meses meses través través reuniones reuniones virtuales, virtuales, definiendo definiendo modificar modificar viejo viejo acuerdo, acuerdo, virtualmente virtualmente suspendido suspendido luego luego cumplieran cumplieran metas metas reservas reservas fiscales fiscales primer primer trimestre trimestre ″, ″, comentaron comentaron expertos expertos research research traders. traders. objetivo objetivo economía economía fondo fondo adelante, adelante, menos, menos, parte parte desembolsos desembolsos comprometidos comprometidos fin fin año año usd usd millones, millones, ayudaría ayudaría reforzar reforzar reservas reservas año año exportaciones exportaciones sector sector agro agro caerían caerían cerca cerca usd usd millones, millones, acuerdo acuerdo proyecciones proyecciones bolsa bolsa comercio comercio rosario rosario bcr. bcr. tales tales desembolsos desembolsos dudas dudas cuánto cuánto podría podría utilizar utilizar intervenir intervenir mercado mercado cambiario. cambiario. fondo fondo aceptaría aceptaría bcra bcra intervenga intervenga eventuales eventuales situaciones situaciones stress. stress. todavía todavía resuelto resuelto monto monto fmi fmi consentiría consentiría intervenciones. intervenciones. reservas reservas internacionales internacionales brutas brutas banco banco central central crecieron crecieron semana semana pasada pasada usd usd millones millones finalizaron finalizaron millones millones dólares. dólares. informe informe anker anker latinoamérica latinoamérica subrayó subrayó cuanto cuanto reservas reservas líquidas líquidas bcra, bcra, “su “su disponibilidad disponibilidad hoy hoy dada dada gran gran medida medida encajes encajes cuentas cuentas bancarias bancarias monedas monedas depositados depositados bcra bcra -usd -usd millones-. millones-. estabilidad estabilidad depósitos depósitos moneda moneda extranjera extranjera crucial crucial sostener sostener capacidad capacidad intervención intervención bcra”. bcra”. luego luego diversas diversas gestiones gestiones ministro ministro massa massa incentivar incentivar pago pago importaciones importaciones divisa divisa china, china, expandió expandió uso uso yuanes yuanes comercio comercio exterior. exterior. enero enero mayo mayo operaciones operaciones equivalente equivalente usd usd millones. millones. además, además, dos dos empresas empresas fabricantes fabricantes electrónica electrónica confirmaron confirmaron pagarán pagarán compromisos compromisos deuda deuda total total usd usd millones millones moneda, moneda, alivio alivio arcas arcas bcra bcra superará superará usd usd millones millones primer primer trimestre trimestre seguir seguir leyendo: leyendo: urls urls imagenes:'

# Row 3 (Ranking 2): only the score changes
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '0.068'

# Row 4 (Ranking 3): score, document name, and content all change
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '0.053'
$ws.Range("C4").Value = 'Noticia N° 06.txt'
$ws.Range("D4").Value = 'título: inflación inflación mayo mayo argentina argentina superó superó venezuela, venezuela, según según estimaciones estimaciones privadas privadas resumen: resumen: abril abril registrado registrado resultado; resultado; variado variado precios precios último último año año ambos ambos países países espera espera diciembre diciembre contenido: contenido: jun, jun, martín martín kanenguiser kanenguiser inflación inflación argentina argentina mayo mayo superó superó venezuela, venezuela, vez vez más, más, según según estimaciones estimaciones privadas privadas ambos ambos países. países. dato dato observatorio observatorio venezolano venezolano finanzas finanzas ovf, ovf, siguen siguen analistas analistas debido debido falta falta rigor rigor estadísticas estadísticas oficiales, oficiales, marcó marcó suba suba precios precios mes mes pasado, pasado, desaceleración desaceleración últimos últimos meses. meses. argentina, argentina, según según mayoría mayoría estimaciones, estimaciones, dato dato mes mes pasado pasado rondó rondó %, %, mientras mientras gobierno gobierno reza reza ubique ubique levemente levemente debajo debajo cifra. cifra. abril, abril, argentina argentina superado superado venezuela: venezuela: registró registró %, %, frente frente informado informado ovf ovf informado informado banco banco central central venezuela, venezuela, daba daba conocer conocer información información octubre octubre año año pasado. pasado. últimos últimos meses meses mayo-mayo, mayo-mayo, según según ovf, ovf, inflación inflación país país gobernado gobernado nicolás nicolás maduro maduro llegó llegó %, %, baja baja respecto respecto registrado registrado abril. abril. además, además, enero enero suba suba precios precios acumulada acumulada %. %. tendencias tendencias inflación inflación últimos últimos meses meses divergentes: divergentes: marzo marzo ovf ovf registró registró venezuela venezuela %, %, abril abril mayo mayo %; %; tanto, tanto, indec indec registró registró argentina argentina marzo, marzo, abril abril estima estima mes mes quedará quedará torno torno interanual; interanual; conocerá conocerá dato dato organismo organismo lidera lidera marco marco lavagna. lavagna. según según ovf, ovf, “en “en medio medio economía economía evidentes evidentes signos signos contracción, contracción, tasa tasa inflación inflación da da tregua. tregua. así, así, mayo mayo tasa tasa inflación inflación mensual mensual alcanzó alcanzó %, %, acumulada acumulada anualizada anualizada %. %. respecto respecto abril, abril, tasa tasa mensual mensual inflación inflación triplicó triplicó aunque aunque tasa tasa doce doce meses meses desaceleró”. desaceleró”. “este “este comportamiento comportamiento inflación inflación ocurriendo ocurriendo entorno entorno signos signos significativa significativa caída caída nivel nivel actividad actividad indiscutibles, indiscutibles, ostensibles ostensibles menores menores ventas ventas comercio comercio producción producción industrial industrial contrajo contrajo primer primer trimestre trimestre ″, ″, indicó. indicó. además, además, “la “la debilidad debilidad demanda demanda agregada agregada notoria notoria debido debido salarios salarios pensiones pensiones sector sector público público pronunciadamente pronunciadamente rezagados rezagados respecto respecto inflación inflación ejecución ejecución gasto gasto parte parte gobierno gobierno baja”. baja”. “al “al comparar comparar cifra cifra inflación inflación mayo mayo respecto respecto aumento aumento tipo tipo cambio cambio mismo mismo mes, mes, clara clara estrecha estrecha relación relación ambas ambas variables, variables, obstante, obstante, mencionó, mencionó, debilidad debilidad demanda, demanda, explicada explicada política política compresión compresión salarial salarial venido venido aplicando aplicando gobierno”, gobierno”, subrayó subrayó ovf. ovf. cuanto cuanto principales principales componentes componentes conforman conforman índice índice nacional nacional precios precios consumidor, consumidor, destacaron destacaron “los “los aumentos aumentos experimentados experimentados rubros rubros esparcimiento esparcimiento %, %, vestido vestido calzado calzado %, %, equipamiento equipamiento hogar hogar alquiler alquiler vivienda vivienda %. %. alimentos alimentos alzas alzas modestas, modestas, mayo mayo incrementaron incrementaron %”. %”. “todo “todo ello ello pone pone manifiesto, manifiesto, política política económica económica aplicada aplicada sido sido ineficaz ineficaz contener contener alza alza precios, precios, aún aún retracción retracción inducida inducida demanda demanda agregada”, agregada”, concluyó concluyó organismo organismo independiente independiente régimen régimen autoritario autoritario maduro. maduro. abril, abril, banco banco central central venezuela venezuela informado informado aumento aumento cuatro cuatro meses, meses, luego luego difundir difundir datos datos medio medio año. año. tanto, tanto, mayoría mayoría relevamientos relevamientos precios precios consumidor consumidor consultoras consultoras argentinas argentinas anticipan anticipan ipc ipc torno torno %. %. c&t c&t indicó indicó relevamiento relevamiento precios precios minoristas minoristas región región gba gba “presentó “presentó alza alza mensual mensual %, %, superando superando largamente largamente variación variación abril abril mayo mayo año año pasado. pasado. así, así, variación variación doce doce meses meses trepó trepó %, %, mayor mayor agosto agosto ″. ″. “el “el rubro rubro mayor mayor incremento incremento vivienda vivienda %, %, reflejando reflejando subas subas gas gas electricidad electricidad principalmente”, principalmente”, aclaró. aclaró. bienes bienes servicios servicios varios varios “ocuparon “ocuparon segundo segundo lugar, lugar, alza alza %, %, explicada explicada cigarrillos cigarrillos artículos artículos tocador”. tocador”. tanto, tanto, “el “el comportamiento comportamiento esparcimiento esparcimiento fuertemente fuertemente influido influido alza alza dólares dólares financieros financieros fines fines abril, abril, vio vio reflejado reflejado turismo turismo productos productos electrónicos”. electrónicos”. vez, vez, “alimentos “alimentos bebidas bebidas creció creció mes. mes. arrancó arrancó mes mes gran gran impulso impulso luego luego moderando. moderando. verduras, verduras, lácteos lácteos derivados derivados harina harina destacaron, destacaron, igual igual alimentos alimentos consumidos consumidos hogar hogar llevar”. llevar”. “en “en salud salud destacó destacó incremento incremento medicamentos, medicamentos, sumó sumó nuevo nuevo ajuste ajuste prepagas”, prepagas”, indicó indicó c&t. c&t. “en “en equipamiento equipamiento mantenimiento mantenimiento hogar, hogar, artefactos artefactos hogar hogar rol rol preponderante preponderante mano mano alza alza dólares dólares financieros financieros pesar pesar liquidaciones liquidaciones mes”, mes”, concluyó. concluyó. parte, parte, ecogo ecogo informó informó si si bien bien última última semana semana mes mes pasado pasado registró registró fuerte fuerte desaceleración desaceleración lugar lugar variación variación alimentos alimentos respecto respecto semana semana anterior, anterior, términos términos generales generales inflación inflación sido sido ciento. ciento. particular, particular, aumento aumento precios precios alimentos alimentos sido sido ciento. ciento. “si “si consideramos consideramos además además evolución evolución alimentos alimentos consumidos consumidos hogar hogar %, %, inflación inflación alimentos alimentos alcanzaría alcanzaría %”, %”, aclaró aclaró consultora consultora dirige dirige marina marina dal dal poggeto. poggeto. lado, lado, lcg lcg detalló detalló “el “el índice índice alimentos alimentos bebidas bebidas presentó presentó inflación inflación mensual mensual promedio promedio últimas últimas semanas semanas punta punta punta punta mismo mismo período”. período”. mes mes sumarán sumarán aumentos aumentos precios precios servicios servicios transporte, transporte, prepagas, prepagas, colegios colegios privados privados combustibles, combustibles, valores valores regulados. regulados. equipo equipo económico económico trata trata controlar controlar cuestión cuestión cambiaria cambiaria inflación inflación acelere acelere todavía todavía motivo motivo resiste resiste pedido pedido fmi fmi acelerar acelerar devaluación devaluación tipo tipo cambio cambio oficial. oficial. seguir seguir leyendo: leyendo: urls urls imagenes:'
